{"js": "// Replace the three-digit x one-digit multiplication answers in the table\n// with the new values, matched by their exact previous text so each cell\n// keeps its own run formatting (font/size) untouched.\nconst replacements = [\n  [\"512\u00d73=1536\", \"633\u00d74=2532\"],\n  [\"221\u00d72=442\", \"325\u00d79=2925\"],\n  [\"206\u00d74=824\", \"176\u00d77=1232\"],\n  [\"543\u00d74=2172\", \"345\u00d78=2760\"],\n  [\"775\u00d78=6200\", \"176\u00d76=1056\"],\n  [\"991\u00d75=4955\", \"854\u00d78=6832\"],\n  [\"730\u00d76=4380\", \"234\u00d75=1170\"],\n  [\"655\u00d74=2620\", \"452\u00d79=4068\"],\n  [\"333\u00d77=2331\", \"423\u00d77=2961\"],\n  [\"648\u00d78=5184\", \"669\u00d74=2676\"],\n  [\"257\u00d77=1799\", \"727\u00d76=4362\"],\n  [\"497\u00d78=3976\", \"188\u00d74=752\"],\n  [\"604\u00d77=4228\", \"668\u00d76=4008\"],\n  [\"739\u00d73=2217\", \"735\u00d73=2205\"],\n  [\"436\u00d74=1744\", \"554\u00d72=1108\"],\n  [\"451\u00d78=3608\", \"713\u00d72=1426\"],\n  [\"986\u00d77=6902\", \"139\u00d72=278\"],\n  [\"659\u00d78=5272\", \"190\u00d79=1710\"],\n  [\"530\u00d76=3180\", \"504\u00d74=2016\"],\n  [\"738\u00d73=2214\", \"866\u00d76=5196\"],\n  [\"346\u00d78=2768\", \"635\u00d79=5715\"],\n  [\"988\u00d76=5928\", \"873\u00d79=7857\"],\n  [\"721\u00d73=2163\", \"552\u00d77=3864\"],\n  [\"608\u00d76=3648\", \"499\u00d76=2994\"],\n  [\"254\u00d76=1524\", \"586\u00d77=4102\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the three-digit x one-digit multiplication answers in the table\n# with the new values, matched by their exact previous text via Find/Replace\n# so each cell keeps its own run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"512\u00d73=1536\", \"633\u00d74=2532\"),\n    @(\"221\u00d72=442\", \"325\u00d79=2925\"),\n    @(\"206\u00d74=824\", \"176\u00d77=1232\"),\n    @(\"543\u00d74=2172\", \"345\u00d78=2760\"),\n    @(\"775\u00d78=6200\", \"176\u00d76=1056\"),\n    @(\"991\u00d75=4955\", \"854\u00d78=6832\"),\n    @(\"730\u00d76=4380\", \"234\u00d75=1170\"),\n    @(\"655\u00d74=2620\", \"452\u00d79=4068\"),\n    @(\"333\u00d77=2331\", \"423\u00d77=2961\"),\n    @(\"648\u00d78=5184\", \"669\u00d74=2676\"),\n    @(\"257\u00d77=1799\", \"727\u00d76=4362\"),\n    @(\"497\u00d78=3976\", \"188\u00d74=752\"),\n    @(\"604\u00d77=4228\", \"668\u00d76=4008\"),\n    @(\"739\u00d73=2217\", \"735\u00d73=2205\"),\n    @(\"436\u00d74=1744\", \"554\u00d72=1108\"),\n    @(\"451\u00d78=3608\", \"713\u00d72=1426\"),\n    @(\"986\u00d77=6902\", \"139\u00d72=278\"),\n    @(\"659\u00d78=5272\", \"190\u00d79=1710\"),\n    @(\"530\u00d76=3180\", \"504\u00d74=2016\"),\n    @(\"738\u00d73=2214\", \"866\u00d76=5196\"),\n    @(\"346\u00d78=2768\", \"635\u00d79=5715\"),\n    @(\"988\u00d76=5928\", \"873\u00d79=7857\"),\n    @(\"721\u00d73=2163\", \"552\u00d77=3864\"),\n    @(\"608\u00d76=3648\", \"499\u00d76=2994\"),\n    @(\"254\u00d76=1524\", \"586\u00d77=4102\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
